$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Logs": append two new mail-log rows (3 and 4)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Klacht over levering"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$logs.Range("D3").Value = "Klacht"
$logs.Range("F3").Value = "2025-06-17 08:58:06"
$logs.Range("G3").Value = "Nee"

$logs.Range("A4").Value = "Offerte voor zakelijke samenwerking"
$logs.Range("B4").Value = "mailmind.test@zohomail.eu"
$logs.Range("C4").Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Range("D4").Value = "Bestelling"
$logs.Range("F4").Value = "2025-06-17 09:28:13"
$logs.Range("G4").Value = "Nee"

# Extend the conditional-formatting ranges from the single data row to
# cover the newly added rows too (D2 -> D2:D4, G2 -> G2:G4).
$catFc = $logs.Range("D2").FormatConditions
for ($i = 1; $i -le $catFc.Count; $i++) {
    $catFc.Item($i).ModifyAppliesToRange($logs.Range("D2:D4"))
}

$answeredFc = $logs.Range("G2").FormatConditions
for ($i = 1; $i -le $answeredFc.Count; $i++) {
    $answeredFc.Item($i).ModifyAppliesToRange($logs.Range("G2:G4"))
}

# ---------------------------------------------------------------------
# Sheet "Dashboard": add the matching category counts
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Klacht"
$dash.Range("B3").Value = 1

$dash.Range("A4").Value = "Bestelling"
$dash.Range("B4").Value = 1

# Grow the bar chart's category/value series references so the chart
# picks up the two new Dashboard rows ($A$2 -> $A$2:$A$4, $B$2 -> $B$2:$B$4).
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$4"
$series.Values = "='Dashboard'!`$B`$2:`$B`$4"
